$d = $word.ActiveDocument

# --- 1) Merge "Public IP(" + _GoBack bookmark + "s)" back into a single
#        run "Public IP(s)".
# The invoice body repeats a "Public IP(s)" table cell several times; only
# one of them still carries a leftover bookmark splitting the text across
# two runs ("Public IP(" / bookmark / "s)"). A document-wide Find/Replace
# for the same literal text is a no-op on every cell that is already a
# single run, and merges the two runs (dropping the stray bookmark) on the
# one cell that still has the split.
$d.Content.Find.Execute("Public IP(s)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Public IP(s)", 2)

# --- 2) Footer: fix the "Grand Coolee Hosting" typo to "Grand Coulee
#        Hosting", and wrap the corrected word "Coulee" in a fresh
#        "_GoBack" bookmark (mirroring Word's last-edit marker). ---
$footer = $d.Sections.Item(1).Footers.Item(1)

$footer.Range.Find.Execute("Coolee", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "Coulee", 2)

$coulee = $footer.Range
if ($coulee.Find.Execute("Coulee")) {
    $d.Bookmarks.Add("_GoBack", $coulee)
}
